$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "LP1912" : append rows 287-303, bump header timestamp/count
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 14:45:41"
$ws1.Cells.Item(3, 1).Value = "Total filas: 302"

$data1 = @(
  @(287, "14:45:30", "14:56", "16_P MOR-SANTA ANA", 11, "LP1912", "30/12/2025"),
  @(288, "14:45:30", "14:58", "215B_EL PATO", 13, "LP1912", "30/12/2025"),
  @(289, "14:45:30", "15:00", "81_EL PELIGRO", 15, "LP1912", "30/12/2025"),
  @(290, "14:45:30", "15:05", "10_OLMOS", 20, "LP1912", "30/12/2025"),
  @(291, "14:45:30", "15:05", "23_HERNANDEZ", 20, "LP1912", "30/12/2025"),
  @(292, "14:45:30", "15:06", "16_SANTA ANA", 21, "LP1912", "30/12/2025"),
  @(293, "14:45:30", "15:20", "15_ABASTO", 35, "LP1912", "30/12/2025"),
  @(294, "14:45:30", "15:21", "26_HERNANDEZ", 36, "LP1912", "30/12/2025"),
  @(295, "14:45:30", "15:32", "84_COLONIA URQUIZA-ESC 49", 47, "LP1912", "30/12/2025"),
  @(296, "14:45:30", "15:42", "10_OLMOS", 57, "LP1912", "30/12/2025"),
  @(297, "14:45:30", "15:45", "14_ABASTO", 60, "LP1912", "30/12/2025"),
  @(298, "14:45:30", "15:51", "23_HERNANDEZ", 66, "LP1912", "30/12/2025"),
  @(299, "14:45:30", "16:01", "10_OLMOS", 76, "LP1912", "30/12/2025"),
  @(300, "14:45:30", "16:02", "11_ETCHEVERRY", 77, "LP1912", "30/12/2025"),
  @(301, "14:45:30", "16:04", "23_HERNANDEZ", 79, "LP1912", "30/12/2025"),
  @(302, "14:45:30", "16:20", "215C_EL PATO", 95, "LP1912", "30/12/2025"),
  @(303, "14:45:30", "16:21", "26_HERNANDEZ", 96, "LP1912", "30/12/2025")
)

foreach ($row in $data1) {
  $r = $row[0]
  $ws1.Cells.Item($r, 2).Value = $row[1]
  $ws1.Cells.Item($r, 3).Value = $row[2]
  $ws1.Cells.Item($r, 4).Value = $row[3]
  $ws1.Cells.Item($r, 5).Value = $row[4]
  $ws1.Cells.Item($r, 6).Value = $row[5]
  $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ------------------------------------------------------------------
# Sheet "LP1912-215" : append rows 23-24, bump header timestamp/count
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 14:45:41"
$ws2.Cells.Item(3, 1).Value = "Total filas: 23"

$data2 = @(
  @(23, "30/12/2025", "14:45:30", "14:58", "215B_EL PATO", 13, "LP1912"),
  @(24, "30/12/2025", "14:45:30", "16:20", "215C_EL PATO", 95, "LP1912")
)

foreach ($row in $data2) {
  $r = $row[0]
  $ws2.Cells.Item($r, 2).Value = $row[1]
  $ws2.Cells.Item($r, 3).Value = $row[2]
  $ws2.Cells.Item($r, 4).Value = $row[3]
  $ws2.Cells.Item($r, 5).Value = $row[4]
  $ws2.Cells.Item($r, 6).Value = $row[5]
  $ws2.Cells.Item($r, 7).Value = $row[6]
}

# ------------------------------------------------------------------
# Sheet "6203-6173" : append rows 44-45, bump header timestamp/count
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 14:45:41"
$ws3.Cells.Item(3, 1).Value = "Total filas: 44"

$data3 = @(
  @(44, "30/12/2025", "14:45:36", "14:53", "215D_LA PLATA", 8, "L6203"),
  @(45, "30/12/2025", "14:45:41", "15:34", "215A_LA PLATA", 49, "L6173")
)

foreach ($row in $data3) {
  $r = $row[0]
  $ws3.Cells.Item($r, 2).Value = $row[1]
  $ws3.Cells.Item($r, 3).Value = $row[2]
  $ws3.Cells.Item($r, 4).Value = $row[3]
  $ws3.Cells.Item($r, 5).Value = $row[4]
  $ws3.Cells.Item($r, 6).Value = $row[5]
  $ws3.Cells.Item($r, 7).Value = $row[6]
}
